$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

$ws.Range("N2").Value = 87.72
$ws.Range("N3").Value = 85.58
$ws.Range("N4").Value = 85.15000000000001
$ws.Range("N5").Value = 86.03
$ws.Range("N6").Value = 82.39
$ws.Range("N7").Value = 89.38
$ws.Range("N8").Value = 97.98999999999999
$ws.Range("N9").Value = 113.2
$ws.Range("N10").Value = 107.7
$ws.Range("N11").Value = 92.91
$ws.Range("N12").Value = 79.56
$ws.Range("N13").Value = 67.45
$ws.Range("N14").Value = 54.6
$ws.Range("N15").Value = 35.3
$ws.Range("N16").Value = 20.96
$ws.Range("N17").Value = 19.7
$ws.Range("N18").Value = 36.78
$ws.Range("N19").Value = 71.87
$ws.Range("N20").Value = 88.27
$ws.Range("N21").Value = 109.29
$ws.Range("N22").Value = 116.83
$ws.Range("N23").Value = 134.94
$ws.Range("N24").Value = 125.71
$ws.Range("N25").Value = 107.86
